# Replace the calculated-column formulas in A4 and A5 with their
# plain, literal results (the user typed over the formulas with the
# static numbers that had already been computed).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openTickets")

$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

# Move / collapse the current selection to A4 (previously E4, with the
# view scrolled so column C was the left-most visible column).
[void]$ws.Range("A4").Select()
